$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-wise updates: (row, newPrice, newVolume, newHour)
# Values are assigned with a leading apostrophe to force Excel to keep
# them as literal text (matching the source workbook's inlineStr cells)
# instead of auto-converting numeric-looking / percent-looking strings.
$updates = @(
    @(2, '326.27', '-0.23%', '3'),
    @(3, '44.32', '-1.31%', '3'),
    @(4, '5.499', '-1.39%', '3'),
    @(5, '0.08010', '-1.17%', '3'),
    @(6, '2.011', '6.25%', '3'),
    @(7, '4.294', '-1.15%', '3'),
    @(8, '2.571', '-5.98%', '3'),
    @(9, '0.9479', '-0.53%', '3'),
    @(10, '0.1150', '-2.09%', '3'),
    @(11, $null, '-3.27%', '3'),
    @(12, '12.15', '41.68%', '3'),
    @(13, '0.09715', '-3.14%', '3'),
    @(14, '0.04568', '8.96%', '3'),
    @(15, $null, '-0.03%', '3'),
    @(16, $null, '-0.68%', '3'),
    @(17, $null, '-4.64%', '3'),
    @(18, '0.005851', '-1.20%', '3'),
    @(19, $null, '-6.32%', '3'),
    @(20, $null, '-0.16%', '3'),
    @(21, '0.1406', '2.17%', '3'),
    @(22, '0.2544', '-3.69%', '3'),
    @(23, '0.001243', '-0.17%', '3'),
    @(24, '0.004303', '-5.59%', '3'),
    @(25, '0.0001188', '-3.81%', '3'),
    @(26, '0.0003742', '-6.65%', '3'),
    @(27, $null, $null, '3'),
    @(28, $null, $null, '3'),
    @(29, $null, $null, '3'),
    @(30, $null, $null, '3'),
    @(31, $null, $null, '3'),
    @(32, $null, $null, '3'),
    @(33, $null, $null, '3'),
    @(34, $null, $null, '3'),
    @(35, $null, $null, '3'),
    @(36, $null, $null, '3'),
    @(37, $null, $null, '3'),
    @(38, '0.02565', '-3.26%', '3'),
    @(39, '0.05518', '-1.16%', '3'),
    @(40, '0.007527', '-2.08%', '3'),
    @(41, '0.1390', '-0.29%', '3'),
    @(42, '0.007596', '-32.99%', '3'),
    @(43, '0.002017', '-2.21%', '3'),
    @(44, '0.008509', '-2.16%', '3'),
    @(45, '0.00007115', '0.08%', '3'),
    @(46, '0.00000000750', '-0.49%', '3'),
    @(47, $null, '0.92%', '3'),
    @(48, '0.004106', '16.39%', '3'),
    @(49, '0.00002099', '-0.49%', '3'),
    @(50, '0.0001999', '-0.49%', '3'),
    @(51, $null, $null, '3')
)

foreach ($u in $updates) {
    $row = $u[0]
    $dVal = $u[1]
    $eVal = $u[2]
    $gVal = $u[3]

    if ($null -ne $dVal) {
        $ws.Cells.Item($row, 4).Value = "'" + $dVal
    }
    if ($null -ne $eVal) {
        $ws.Cells.Item($row, 5).Value = "'" + $eVal
    }
    if ($null -ne $gVal) {
        $ws.Cells.Item($row, 7).Value = "'" + $gVal
    }
}

Write-Host "Updated $($updates.Count) rows"
